$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.27007749937006
$ws.Range("D2").Value = 5.43959068619453
$ws.Range("E2").Value = 15.30061576191037
$ws.Range("F2").Value = 28.09751633851083
$ws.Range("G2").Value = 35.54963467961222
$ws.Range("H2").Value = 15.80924733870064
$ws.Range("I2").Value = 27.32656388733457
$ws.Range("K2").Value = 9.446741730515994
$ws.Range("L2").Value = 9.495228469461653
$ws.Range("M2").Value = 14.4636480358399
$ws.Range("N2").Value = 20.58930334059145
$ws.Range("B3").Value = 13.17904548488857
$ws.Range("D3").Value = 5.42302776705454
$ws.Range("E3").Value = 15.30015968497214
$ws.Range("F3").Value = 28.042605111508
$ws.Range("G3").Value = 35.42675836595008
$ws.Range("H3").Value = 15.83900580367749
$ws.Range("I3").Value = 27.41760784283329
$ws.Range("K3").Value = 9.127488161034385
$ws.Range("L3").Value = 9.48275529099222
$ws.Range("M3").Value = 14.44711818294754
$ws.Range("N3").Value = 20.65246199561317
$ws.Range("B4").Value = 13.12531296494148
$ws.Range("D4").Value = 5.412648226768514
$ws.Range("E4").Value = 15.30019364599664
$ws.Range("F4").Value = 28.01655210129555
$ws.Range("G4").Value = 35.36245459578806
$ws.Range("H4").Value = 15.86021439990216
$ws.Range("I4").Value = 27.47763959583413
$ws.Range("K4").Value = 8.927675279479566
$ws.Range("L4").Value = 9.476708544936375
$ws.Range("M4").Value = 14.43913583484856
$ws.Range("N4").Value = 20.69304869166673
$ws.Range("B5").Value = 13.10397987605377
$ws.Range("D5").Value = 5.408366444556976
$ws.Range("E5").Value = 15.30028652501217
$ws.Range("F5").Value = 28.00786905218356
$ws.Range("G5").Value = 35.33907054135451
$ws.Range("H5").Value = 15.86959458910538
$ws.Range("I5").Value = 27.50314169311203
$ws.Range("K5").Value = 8.845432601787573
$ws.Range("L5").Value = 9.474651900212455
$ws.Range("M5").Value = 14.43643109104593
$ws.Range("N5").Value = 20.71004380766935
$ws.Range("B6").Value = 13.10047211947209
$ws.Range("D6").Value = 5.407652336427775
$ws.Range("E6").Value = 15.30030672377015
$ws.Range("F6").Value = 28.00654419599878
$ws.Range("G6").Value = 35.335358429595
$ws.Range("H6").Value = 15.87119667587107
$ws.Range("I6").Value = 27.507439031882
$ws.Range("K6").Value = 8.831731087980286
$ws.Range("L6").Value = 9.474335062729574
$ws.Range("M6").Value = 14.43601516385288
$ws.Range("N6").Value = 20.71289339609151
$ws.Range("B7").Value = 13.12502295294149
$ws.Range("D7").Value = 5.412590690889352
$ws.Range("E7").Value = 15.30019457843074
$ws.Range("F7").Value = 28.01642716131064
$ws.Range("G7").Value = 35.3621277900019
$ws.Range("H7").Value = 15.86033791963352
$ws.Range("I7").Value = 27.47797932043549
$ws.Range("K7").Value = 8.926569247727098
$ws.Range("L7").Value = 9.476679155904066
$ws.Range("M7").Value = 14.43909713436512
$ws.Range("N7").Value = 20.69327604686546
$ws.Range("B8").Value = 13.23825319623895
$ws.Range("D8").Value = 5.43392378014583
$ws.Range("E8").Value = 15.30039338531003
$ws.Range("F8").Value = 28.07699729577509
$ws.Range("G8").Value = 35.50496643843277
$ws.Range("H8").Value = 15.81889798432252
$ws.Range("I8").Value = 27.35709860867432
$ws.Range("K8").Value = 9.337527639350645
$ws.Range("L8").Value = 9.490594541099522
$ws.Range("M8").Value = 14.45750077913298
$ws.Range("N8").Value = 20.61070619278671
$ws.Range("B9").Value = 13.47653774425501
$ws.Range("D9").Value = 5.474061875705432
$ws.Range("E9").Value = 15.30327179037456
$ws.Range("F9").Value = 28.25618255070179
$ws.Range("G9").Value = 35.87249432200937
$ws.Range("H9").Value = 15.76097565790786
$ws.Range("I9").Value = 27.15282244638216
$ws.Range("K9").Value = 10.25409449207247
$ws.Range("L9").Value = 9.53056753431842
$ws.Range("M9").Value = 14.51063296783587
$ws.Range("N9").Value = 20.46306369790679
$ws.Range("B10").Value = 13.66023746070868
$ws.Range("D10").Value = 5.502477351030169
$ws.Range("E10").Value = 15.30690099192561
$ws.Range("F10").Value = 28.42400168437491
$ws.Range("G10").Value = 36.19422110189763
$ws.Range("H10").Value = 15.73269747847825
$ws.Range("I10").Value = 27.02271524483002
$ws.Range("K10").Value = 10.94859024031854
$ws.Range("L10").Value = 9.567514526833966
$ws.Range("M10").Value = 14.55983926363179
$ws.Range("N10").Value = 20.36320694311739
$ws.Range("B11").Value = 13.74539965680053
$ws.Range("D11").Value = 5.515162635304972
$ws.Range("E11").Value = 15.30888037793613
$ws.Range("F11").Value = 28.50802084587039
$ws.Range("G11").Value = 36.35139303580958
$ws.Range("H11").Value = 15.72294072784665
$ws.Range("I11").Value = 26.96786313277379
$ws.Range("K11").Value = 11.24854422335192
$ws.Range("L11").Value = 9.585930245368358
$ws.Range("M11").Value = 14.58438001589687
$ws.Range("N11").Value = 20.31963180264531
$ws.Range("B12").Value = 13.77785359701633
$ws.Range("D12").Value = 5.519930805804231
$ws.Range("E12").Value = 15.30967708758438
$ws.Range("F12").Value = 28.54092243815357
$ws.Range("G12").Value = 36.41242496746644
$ws.Range("H12").Value = 15.71969324270163
$ws.Range("I12").Value = 26.94771545073555
$ws.Range("K12").Value = 11.35980658969318
$ws.Range("L12").Value = 9.593131423020713
$ws.Range("M12").Value = 14.59397809251569
$ws.Range("N12").Value = 20.30339571533382
$ws.Range("B13").Value = 13.77085534425129
$ws.Range("D13").Value = 5.518905485409712
$ws.Range("E13").Value = 15.30950340562246
$ws.Range("F13").Value = 28.53378853614591
$ws.Range("G13").Value = 36.39921398698873
$ws.Range("H13").Value = 15.72037275367843
$ws.Range("I13").Value = 26.95202687295497
$ws.Range("K13").Value = 11.33594795604631
$ws.Range("L13").Value = 9.59157046250005
$ws.Range("M13").Value = 14.59189749284114
$ws.Range("N13").Value = 20.30688068567768
$ws.Range("B14").Value = 13.74806569837232
$ws.Range("D14").Value = 5.515555631080515
$ws.Range("E14").Value = 15.30894497898201
$ws.Range("F14").Value = 28.51070604072535
$ws.Range("G14").Value = 36.35638406688258
$ws.Range("H14").Value = 15.72266459142082
$ws.Range("I14").Value = 26.96619307043715
$ws.Range("K14").Value = 11.25774454670604
$ws.Range("L14").Value = 9.58651815250656
$ws.Range("M14").Value = 14.58516357239002
$ws.Range("N14").Value = 20.31829074826787
$ws.Range("B15").Value = 13.73413232523757
$ws.Range("D15").Value = 5.513499105220942
$ws.Range("E15").Value = 15.30860906595949
$ws.Range("F15").Value = 28.49670810709295
$ws.Range("G15").Value = 36.33034544019117
$ws.Range("H15").Value = 15.72412665347867
$ws.Range("I15").Value = 26.9749515002781
$ws.Range("K15").Value = 11.20953939978321
$ws.Range("L15").Value = 9.583452984808076
$ws.Range("M15").Value = 14.5810784193409
$ws.Range("N15").Value = 20.32531420036828
$ws.Range("B16").Value = 13.65470192834485
$ws.Range("D16").Value = 5.501643418106688
$ws.Range("E16").Value = 15.30677823747858
$ws.Range("F16").Value = 28.41866377205676
$ws.Range("G16").Value = 36.18416376745913
$ws.Range("H16").Value = 15.7333976648731
$ws.Range("I16").Value = 27.0263872065821
$ws.Range("K16").Value = 10.92866368009634
$ws.Range("L16").Value = 9.566343066826748
$ws.Range("M16").Value = 14.55827843518862
$ws.Range("N16").Value = 20.36609181031111
$ws.Range("B17").Value = 13.60636532403869
$ws.Range("D17").Value = 5.494308022390447
$ws.Range("E17").Value = 15.30573915546984
$ws.Range("F17").Value = 28.37273979154935
$ws.Range("G17").Value = 36.09722827230485
$ws.Range("H17").Value = 15.73988126405248
$ws.Range("I17").Value = 27.05905171299349
$ws.Range("K17").Value = 10.75224310415815
$ws.Range("L17").Value = 9.556256044495225
$ws.Range("M17").Value = 14.54484020422624
$ws.Range("N17").Value = 20.3915805876786
$ws.Range("B18").Value = 13.57871513340567
$ws.Range("D18").Value = 5.490066255468569
$ws.Range("E18").Value = 15.30517242199054
$ws.Range("F18").Value = 28.34704940084515
$ws.Range("G18").Value = 36.04824607800017
$ws.Range("H18").Value = 15.74390290404007
$ws.Range("I18").Value = 27.07824741623092
$ws.Range("K18").Value = 10.64926812886265
$ws.Range("L18").Value = 9.55060588048795
$ws.Range("M18").Value = 14.53731419766854
$ws.Range("N18").Value = 20.40641527766406
$ws.Range("B19").Value = 13.56938006403986
$ws.Range("D19").Value = 5.488626202038876
$ws.Range("E19").Value = 15.30498584827466
$ws.Range("F19").Value = 28.33847593875656
$ws.Range("G19").Value = 36.03183807218874
$ws.Range("H19").Value = 15.74531477577837
$ws.Range("I19").Value = 27.08481680564621
$ws.Range("K19").Value = 10.61414541259322
$ws.Range("L19").Value = 9.548718983404401
$ws.Range("M19").Value = 14.53480108981101
$ws.Range("N19").Value = 20.41146800741357
$ws.Range("B20").Value = 13.61149531093502
$ws.Range("D20").Value = 5.495091236924528
$ws.Range("E20").Value = 15.30584656805929
$ws.Range("F20").Value = 28.37755368917336
$ws.Range("G20").Value = 36.1063773424294
$ws.Range("H20").Value = 15.73916080357413
$ws.Range("I20").Value = 27.05553229741373
$ws.Range("K20").Value = 10.77117909114719
$ws.Range("L20").Value = 9.557314158942983
$ws.Range("M20").Value = 14.54624972111106
$ws.Range("N20").Value = 20.38884924135468
$ws.Range("B21").Value = 13.75475420764071
$ws.Range("D21").Value = 5.516540533607901
$ws.Range("E21").Value = 15.30910772310149
$ws.Range("F21").Value = 28.5174566304394
$ws.Range("G21").Value = 36.36892349750349
$ws.Range("H21").Value = 15.7219792844
$ws.Range("I21").Value = 26.96201518700319
$ws.Range("K21").Value = 11.2807780305931
$ws.Range("L21").Value = 9.587995993888777
$ws.Range("M21").Value = 14.5871332516642
$ws.Range("N21").Value = 20.31493215710882
$ws.Range("B22").Value = 13.84956470143788
$ws.Range("D22").Value = 5.530351761388352
$ws.Range("E22").Value = 15.31151392328115
$ws.Range("F22").Value = 28.61520767624763
$ws.Range("G22").Value = 36.54931801321141
$ws.Range("H22").Value = 15.71335676666571
$ws.Range("I22").Value = 26.90453147389007
$ws.Range("K22").Value = 11.60027768701441
$ws.Range("L22").Value = 9.609372748939512
$ws.Range("M22").Value = 14.61562837177837
$ws.Range("N22").Value = 20.26816640899419
$ws.Range("B23").Value = 13.79886256254884
$ws.Range("D23").Value = 5.522999667738691
$ws.Range("E23").Value = 15.3102045637154
$ws.Range("F23").Value = 28.56246478012383
$ws.Range("G23").Value = 36.45224675835118
$ws.Range("H23").Value = 15.71772017770643
$ws.Range("I23").Value = 26.93487887787787
$ws.Range("K23").Value = 11.43100208051819
$ws.Range("L23").Value = 9.597843685712157
$ws.Range("M23").Value = 14.60025929631091
$ws.Range("N23").Value = 20.29298534139274
$ws.Range("B24").Value = 13.60917560849369
$ws.Range("D24").Value = 5.494737222028353
$ws.Range("E24").Value = 15.30579791140058
$ws.Range("F24").Value = 28.37537510704665
$ws.Range("G24").Value = 36.10223793542484
$ws.Range("H24").Value = 15.73948560754956
$ws.Range("I24").Value = 27.05712212800459
$ws.Range("K24").Value = 10.76262295248259
$ws.Range("L24").Value = 9.556835320997903
$ws.Range("M24").Value = 14.54561185574936
$ws.Range("N24").Value = 20.39008351956787
$ws.Range("B25").Value = 13.41046527909914
$ws.Range("D25").Value = 5.463389159186537
$ws.Range("E25").Value = 15.3022269632886
$ws.Range("F25").Value = 28.2013007996029
$ws.Range("G25").Value = 35.76385831501894
$ws.Range("H25").Value = 15.77414083957941
$ws.Range("I25").Value = 27.20457666519207
$ws.Range("K25").Value = 9.984385109628617
$ws.Range("L25").Value = 9.518410660403719
$ws.Range("M25").Value = 14.49445623831486
$ws.Range("N25").Value = 20.50148554638109
